# 1 add meetvideoview test activity 2 fix multi-seek event just do once problem
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bug")

# --- Fix up formatting for row 34 (G column) and new row 35 ---
# Row 35 inherits the same format pattern used by similar closed-bug rows
# (B/C/D/E/F/G/H = styles 6,6,6,14,6,19,6) -- borrow from row 18, then fix G.
$ws.Range("B18:H18").Copy() | Out-Null
$ws.Range("B35:H35").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("G33").Copy() | Out-Null
$ws.Range("G35").PasteSpecial(-4122) | Out-Null       # xlPasteFormats

# G34's status moves from "tracking" styling (19) to a plain "fixed" styling (15)
$ws.Range("H14").Copy() | Out-Null
$ws.Range("G34").PasteSpecial(-4122) | Out-Null       # xlPasteFormats

$excel.CutCopyMode = 0

# --- Row 34: bug #33 is now resolved; fold the old "stayAwake" note into F,
#     record the new fix status/description ---
$ws.Range("F34").Value = "播放中自动黑屏休眠了(系统播放器还没有解决)"
$ws.Range("G34").Value = "fixed"
$ws.Range("H34").Value = "没有在start stop等 设置stayAwake, 增加了nativeplayer基类"

# --- Row 35: new bug #34 ---
$ws.Range("B35").Value = 34
$ws.Range("C35").Value = "PPBOX-3170"
$ws.Range("D35").Value = "OTT"
$ws.Range("E35").Value = 20150225
$ws.Range("F35").Value = "rmvb片源快进后声音异常"
$ws.Range("G35").Value = "tracking"
$ws.Range("H35").Value = "开始播放后，seek到片子中段，会有杂音"

# --- Update the saved sheet view / selection state ---
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("H36").Select() | Out-Null
